$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 65 (shifts existing rows 65-133 down to 66-134)
$ws.Rows(65).Insert()

# Populate the newly inserted row 65 with the new weekly data point
$ws.Range("A65").Value = 11
$ws.Range("B65").Value = "Vega Monumental Concepción"
$ws.Range("C65").Value = "Bíobío"
$ws.Range("D65").Value = 44586
$ws.Range("E65").Value = 8
$ws.Range("F65").Value = 100112003
$ws.Range("G65").Value = "Ajo"
$ws.Range("H65").Value = "Chino"
$ws.Range("I65").Value = "1a (cosecha)"
$ws.Range("J65").Value = 220
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 16000
$ws.Range("M65").Value = 15545
$ws.Range("N65").Value = "`$/caja 10 kilos"
$ws.Range("O65").Value = "Provincia de Talagante"
$ws.Range("P65").Value = 1554
$ws.Range("Q65").Value = 10
$ws.Range("R65").Value = "Hortaliza"
